$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume 1h) updates per source diff.
# D-column numeric-looking strings are forced to Text format so Excel
# does not reinterpret them as numbers (which would drop formatting such
# as trailing zeros, e.g. "1.00" -> 1, "69.00" -> 69, "0.0860" -> 0.086).

$ws.Range("D2").Value = '42.065.98'
$ws.Range("E2").Value = '  -3.71%  '
$ws.Range("D3").Value = '2.195.39'
$ws.Range("E3").Value = '  -3.72%  '
$ws.Range("E4").Value = '  +0.23%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '105.97'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -14.50%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '291.54'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +9.59%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.618'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -3.11%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("E9").Value = '  -6.33%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '43.51'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -9.59%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0902'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -4.63%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '53.98'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.61%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '8.63'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -6.83%  '
$ws.Range("E14").Value = '  -3.70%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.921'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.29%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '14.71'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -4.65%  '
$ws.Range("D17").Value = '2.525.72'
$ws.Range("E17").Value = '  -3.74%  '
$ws.Range("D18").Value = '2.207.56'
$ws.Range("E18").Value = '  -3.24%  '
$ws.Range("D19").Value = '41.960.09'
$ws.Range("E19").Value = '  -3.89%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.15'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("E21").Value = '  -5.87%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '72.12'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -0.29%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '3.38'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +17.35%  '
$ws.Range("E24").Value = '  -8.19%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '225.65'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -4.28%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.84'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -6.06%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("E28").Value = '  -3.47%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '3.88'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.18%  '
$ws.Range("E30").Value = '  -1.67%  '
$ws.Range("E31").Value = '  -5.03%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '36.98'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -13.53%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '170.93'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -1.02%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '20.65'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.76%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.0860'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.93%  '
$ws.Range("E36").Value = '  -5.83%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.89'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +5.01%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.16'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.49%  '
$ws.Range("E39").Value = '  -4.27%  '
$ws.Range("E40").Value = '  -5.63%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.101'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -5.68%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -5.14%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '69.00'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -7.90%  '
$ws.Range("E44").Value = '  -5.53%  '
$ws.Range("E45").Value = '  +0.17%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '12.52'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -10.07%  '
$ws.Range("E47").Value = '  -6.90%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '5.34'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.56%  '
$ws.Range("E49").Value = '  +1.82%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '101.44'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -0.47%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '8.28'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.67%  '
